$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.924.83"
$ws.Range("E2").Value = "  +4.95%  "

Set-TextValue $ws.Range("D3") "2.276.89"
$ws.Range("E3").Value = "  +2.66%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue $ws.Range("D5") "303.84"
$ws.Range("E5").Value = "  +3.57%  "

Set-TextValue $ws.Range("D6") "93.44"
$ws.Range("E6").Value = "  +7.92%  "

Set-TextValue $ws.Range("D7") "0.531"
$ws.Range("E7").Value = "  +4.05%  "

$ws.Range("E8").Value = "  -0.03%  "

Set-TextValue $ws.Range("D9") "0.487"
$ws.Range("E9").Value = "  +4.48%  "

Set-TextValue $ws.Range("D10") "32.57"
$ws.Range("E10").Value = "  +6.97%  "

Set-TextValue $ws.Range("D11") "53.18"
$ws.Range("E11").Value = "  +5.95%  "

Set-TextValue $ws.Range("D12") "0.0803"
$ws.Range("E12").Value = "  +2.84%  "

Set-TextValue $ws.Range("D13") "0.116"
$ws.Range("E13").Value = "  +2.87%  "

Set-TextValue $ws.Range("D14") "6.71"
$ws.Range("E14").Value = "  +4.71%  "

Set-TextValue $ws.Range("D15") "2.628.83"
$ws.Range("E15").Value = "  +2.73%  "

Set-TextValue $ws.Range("D16") "14.30"
$ws.Range("E16").Value = "  +3.66%  "

Set-TextValue $ws.Range("D17") "2.283.50"
$ws.Range("E17").Value = "  +2.43%  "

$ws.Range("E18").Value = "  +3.77%  "

Set-TextValue $ws.Range("D19") "41.828.62"
$ws.Range("E19").Value = "  +4.93%  "

Set-TextValue $ws.Range("D20") "12.32"
$ws.Range("E20").Value = "  +9.91%  "

Set-TextValue $ws.Range("D21") "0.0₃0908"
$ws.Range("E21").Value = "  +2.79%  "

Set-TextValue $ws.Range("D22") "5.97"
$ws.Range("E22").Value = "  +3.83%  "

Set-TextValue $ws.Range("D23") "67.44"
$ws.Range("E23").Value = "  +2.81%  "

Set-TextValue $ws.Range("D24") "244.46"
$ws.Range("E24").Value = "  +3.65%  "

Set-TextValue $ws.Range("D25") "2.58"
$ws.Range("E25").Value = "  +4.89%  "

Set-TextValue $ws.Range("D26") "1.94"
$ws.Range("E26").Value = "  +6.06%  "

$ws.Range("E27").Value = "  -0.06%  "

Set-TextValue $ws.Range("D28") "24.40"
$ws.Range("E28").Value = "  +6.04%  "

Set-TextValue $ws.Range("D29") "9.65"
$ws.Range("E29").Value = "  +4.33%  "

$ws.Range("E30").Value = "  -10.71%  "

# Rows 31 and 32 swap coins (Monero <-> InjectiveProtocol) along with updated values
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D31") "34.23"
$ws.Range("E31").Value = "  +8.60%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D32") "158.51"
$ws.Range("E32").Value = "  +0.38%  "

Set-TextValue $ws.Range("D33") "0.999"
$ws.Range("E33").Value = "  +0.01%  "

Set-TextValue $ws.Range("D34") "5.21"
$ws.Range("E34").Value = "  +5.34%  "

Set-TextValue $ws.Range("D35") "0.0751"
$ws.Range("E35").Value = "  +5.84%  "

Set-TextValue $ws.Range("D36") "3.08"
$ws.Range("E36").Value = "  +2.10%  "

$ws.Range("E37").Value = "  +2.64%  "

Set-TextValue $ws.Range("D38") "16.82"
$ws.Range("E38").Value = "  +9.45%  "

Set-TextValue $ws.Range("D39") "0.105"
$ws.Range("E39").Value = "  +6.23%  "

Set-TextValue $ws.Range("D40") "0.116"
$ws.Range("E40").Value = "  +3.19%  "

Set-TextValue $ws.Range("D41") "1.83"
$ws.Range("E41").Value = "  +5.87%  "

Set-TextValue $ws.Range("D42") "3.95"
$ws.Range("E42").Value = "  +6.79%  "

Set-TextValue $ws.Range("D43") "2.072.19"
$ws.Range("E43").Value = "  -0.56%  "

Set-TextValue $ws.Range("D44") "19.72"
$ws.Range("E44").Value = "  +10.86%  "

Set-TextValue $ws.Range("D45") "0.0281"
$ws.Range("E45").Value = "  +4.33%  "

Set-TextValue $ws.Range("D46") "10.40"
$ws.Range("E46").Value = "  +4.37%  "

$ws.Range("E47").Value = "  +9.18%  "

$ws.Range("E48").Value = "  +3.51%  "

Set-TextValue $ws.Range("D49") "73.63"
$ws.Range("E49").Value = "  +9.50%  "

Set-TextValue $ws.Range("D50") "1.54"
$ws.Range("E50").Value = "  +5.57%  "

$ws.Range("E51").Value = "  +3.63%  "
